$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Trim the "Descrição" paragraph: drop the trailing clause about the
#    database, keeping the final "." run untouched.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Descrição:  O Software tem o intuito de ajudar a gerenciar uma escola, desde os funcionários até os alunos, contando com um banco de dados para armazenar os dados",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Descrição:  O Software tem o intuito de ajudar a gerenciar uma escola, desde os funcionários até os alunos",
    2) | Out-Null

# ------------------------------------------------------------------
# 2. Rewrite the requirement bullet paragraphs (paragraphs 4-9) with
#    their new wording.
# ------------------------------------------------------------------
$d.Paragraphs(4).Range.Text = "- Apresentar a capacidade de realizar o registro, edição e deleção dos funcionários, professores e alunos."
$d.Paragraphs(5).Range.Text = "- Gerenciar a relação das disciplinas, bem como os horários e a grade curricular da escola."
$d.Paragraphs(6).Range.Text = "- Gerenciar as faltas, suspensão e histórico escolar de um aluno."
$d.Paragraphs(7).Range.Text = "- Gerar um CALENDÁRIO escolar"
$d.Paragraphs(8).Range.Text = "- O sistema deve ser capaz de suspender alunos"
$d.Paragraphs(9).Range.Text = "- Gerenciar salários, atrasos e pagamentos."

# ------------------------------------------------------------------
# 3. Append a brand-new requirement paragraph at the end of the
#    document body (after the last bullet, before the sectPr).
# ------------------------------------------------------------------
$endOfBody = $d.Paragraphs(9).Range
$endOfBody.Collapse(0)
$endOfBody.InsertParagraphAfter()
$newBullet = $d.Paragraphs(10).Range
$newBullet.Text = "- Deve apresentar opções de recuperação de senha."

# ------------------------------------------------------------------
# 4. Materialize the built-in "List Paragraph" style definition into
#    styles.xml (it ends up unused by any paragraph, matching the
#    source edit which only touched the style catalogue).
# ------------------------------------------------------------------
$style = $d.Styles.Add("PargrafodaLista", 1)
$style.NameLocal = "List Paragraph"
$style.BaseStyle = "Normal"
$style.Priority = 34
$style.QuickStyle = $true
$style.ParagraphFormat.LeftIndent = 36

$tailRange = $d.Content
$tailRange.Collapse(0)
$tailRange.InsertParagraphAfter()
$scratchPara = $d.Paragraphs($d.Paragraphs.Count)
$scratchPara.Style = $style
$scratchPara.Range.Delete()
